# Build site at 2022-09-26 16:07:08 UTC
# Applies the LOT2008.xlsx content update:
#  - Removes the detailed "Programa resumido" paragraph row (old row 14),
#    which shifts every following row up by one.
#  - Updates several cells to their new values as a consequence of the
#    content re-shuffle (teacher name / "Semestral" / reordered evaluation
#    paragraphs, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Teacher name now appears right under "Objetivos:" (row 10)
$ws.Cells.Item(10, 2).Value = "5840494 - Maria Eleonora Andrade de Carvalho"
$ws.Cells.Item(10, 3).Value = "5840494 - Maria Eleonora Andrade de Carvalho"

# 2) Delete the old detailed "Programa resumido" paragraph row (row 14).
#    This shifts all subsequent rows up by one.
$ws.Rows.Item(14).Delete()

# 3) Row 13 becomes "Programa resumido:" / "Semestral"
$ws.Cells.Item(13, 1).Value = "Programa resumido:"
$ws.Cells.Item(13, 2).Value = "Semestral"
$ws.Cells.Item(13, 3).Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# 4) Row 15 ("Programa:") now shows the date value instead of the long text
$ws.Cells.Item(15, 2).Value = "01/01/2018"
$ws.Cells.Item(15, 3).Value = "01/01/2018"

# 5) Evaluation block reshuffle (rows 18-21)
$ws.Cells.Item(18, 2).Value = "5840494 - Maria Eleonora Andrade de Carvalho"
$ws.Cells.Item(18, 3).Value = "5840494 - Maria Eleonora Andrade de Carvalho"

$ws.Cells.Item(19, 2).Value = "A avaliação será feita por meio de provas escritas."
$ws.Cells.Item(19, 3).Value = "A avaliação será feita por meio de provas escritas."

$ws.Cells.Item(20, 2).Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2"
$ws.Cells.Item(20, 3).Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2"

$ws.Cells.Item(21, 2).Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Cells.Item(21, 3).Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
